$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Essence of Ebonhawke (name unchanged)
$ws.Range("C2").Value = 1900004
$ws.Range("D2").Value = 4349998
$ws.Range("E2").Value = 179.74943
$ws.Range("F2").Value = 94.60476398997054
$ws.Range("G2").Value = "OK_SLOW"
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 0.323
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.6204853842040662

# Row 3: Shrouded Bench of the Final Judge -> Glyph of the Unbound (Unused)
$ws.Range("A3").Value = "Glyph of the Unbound (Unused)"
$ws.Range("B3").Value = 25441
$ws.Range("C3").Value = 2010002
$ws.Range("D3").Value = 2888483
$ws.Range("E3").Value = 44.52085499999998
$ws.Range("F3").Value = 22.14965706501783
$ws.Range("J3").Value = $false

# Row 4: Celestial Infusion (Blue) -> Ignus Fatuus
$ws.Range("A4").Value = "Ignus Fatuus"
$ws.Range("B4").Value = 15717
$ws.Range("C4").Value = 2696969
$ws.Range("D4").Value = 3669998
$ws.Range("E4").Value = 42.25292999999998
$ws.Range("F4").Value = 15.66682079030199
$ws.Range("G4").Value = "OK_LIQUID"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0.9360000000000001
$ws.Range("L4").Value = 0.9396660624621362
$ws.Range("M4").Value = ""

# Row 5: Crustacea -> Fox Fire Spear Skin
$ws.Range("A5").Value = "Fox Fire Spear Skin"
$ws.Range("B5").Value = 332607
$ws.Range("C5").Value = 850010
$ws.Range("D5").Value = 1449897
$ws.Range("E5").Value = 38.24024499999999
$ws.Range("F5").Value = 44.98799425889107
$ws.Range("G5").Value = "HOLD_RISK"
$ws.Range("H5").Value = 98
$ws.Range("I5").Value = 0.01
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.03036139133931093

# Row 6: Advanced Leather Rack -> Arthropoda
$ws.Range("A6").Value = "Arthropoda"
$ws.Range("B6").Value = 19936
$ws.Range("C6").Value = 965729
$ws.Range("D6").Value = 1366670
$ws.Range("E6").Value = 19.59405
$ws.Range("F6").Value = 20.28938760252617
$ws.Range("I6").Value = 1.097
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.9627948221649534

# Row 7: Dragon''s Claw -> Winter''s Heart Infusion
$ws.Range("A7").Value = "Winter''s Heart Infusion"
$ws.Range("B7").Value = 24511
$ws.Range("C7").Value = 1657076
$ws.Range("D7").Value = 2149999
$ws.Range("E7").Value = 17.04231499999999
$ws.Range("F7").Value = 10.28457053267321
$ws.Range("I7").Value = 2.385
$ws.Range("L7").Value = 0.9992193330446831

# Row 8: Collapsing Star Spear Skin -> +14 Agony Infusion
$ws.Range("A8").Value = "+14 Agony Infusion"
$ws.Range("B8").Value = 21827
$ws.Range("C8").Value = 1510004
$ws.Range("D8").Value = 1899999
$ws.Range("E8").Value = 10.49951499999999
$ws.Range("F8").Value = 6.953302772707881
$ws.Range("I8").Value = 2.106
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.9981970142173253

# Row 9: Vintage Black Lion Weapon Box -> Elder Wood Logging Node
$ws.Range("A9").Value = "Elder Wood Logging Node"
$ws.Range("B9").Value = 22586
$ws.Range("C9").Value = 1000003
$ws.Range("D9").Value = 1299750
$ws.Range("E9").Value = 10.47845
$ws.Range("F9").Value = 10.47841856474431
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 0.733
$ws.Range("L9").Value = 0.889031248563259

# Row 10: Ghostly Infusion (name unchanged)
$ws.Range("B10").Value = 24213
$ws.Range("C10").Value = 609994
$ws.Range("D10").Value = 839990
$ws.Range("E10").Value = 10.39975
$ws.Range("F10").Value = 17.04893818627724
$ws.Range("I10").Value = 8.676
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0.9999999999950331

# Row 11: Zojja''s Berserker Insignia -> Winter''s Heart Infusion
$ws.Range("A11").Value = "Winter''s Heart Infusion"
$ws.Range("B11").Value = 24517
$ws.Range("C11").Value = 880725
$ws.Range("D11").Value = 1158090
$ws.Range("E11").Value = 10.36515
$ws.Range("F11").Value = 11.76888359022396
$ws.Range("I11").Value = 2.965
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.9998629051024138
